# NYPD CompStat weekly update — "New crime data collected"
# Updates the report header (volume/week-ending text) and refreshes every
# crime-statistic figure in the table (rows 15-31) to the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: volume/issue number and the report date range.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# ---------------------------------------------------------------------
# Helper: turn a currently-blank("0"/"***.*" placeholder text) cell into
# a real number while keeping the workbook's normal numeric look
# (reuses the number format of a same-column/style neighbour cell).
# ---------------------------------------------------------------------
function Set-NumberCell($a1, $value, $formatSourceA1) {
    $ws.Range($formatSourceA1).Copy() | Out-Null
    $ws.Range($a1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($a1).Value = $value
}

# ---------------------------------------------------------------------
# Helper: turn a currently-numeric cell back into the blank placeholder
# text ("0" or "***.*") used throughout the sheet for n/a figures.
# ---------------------------------------------------------------------
function Set-PlaceholderCell($a1, $text, $formatSourceA1) {
    $ws.Range($a1).Value = "'" + $text
    $ws.Range($formatSourceA1).Copy() | Out-Null
    $ws.Range($a1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------
# Row 15 — Rape
# ---------------------------------------------------------------------
Set-NumberCell "C15" 2 "J15"
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 28
$ws.Range("K15").Value = 47.368421052631
$ws.Range("L15").Value = 180
$ws.Range("M15").Value = 211.111111111111
$ws.Range("N15").Value = -9.677419354838

# ---------------------------------------------------------------------
# Row 16 — Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 14.285714285714
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = -26.470588235294
$ws.Range("I16").Value = 229
$ws.Range("J16").Value = 282
$ws.Range("K16").Value = -18.794326241134
$ws.Range("L16").Value = -30.815709969788
$ws.Range("M16").Value = 110.091743119266
$ws.Range("N16").Value = -87.192393736017

# ---------------------------------------------------------------------
# Row 17 — Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 14
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 250
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 41.176470588235
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 373
$ws.Range("K17").Value = 7.238605898123
$ws.Range("L17").Value = 14.613180515759
$ws.Range("M17").Value = 194.117647058824
$ws.Range("N17").Value = -13.232104121475

# ---------------------------------------------------------------------
# Row 18 — Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 42
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 299
$ws.Range("J18").Value = 260
$ws.Range("K18").Value = 15
$ws.Range("L18").Value = 4.181184668989
$ws.Range("M18").Value = 33.482142857142
$ws.Range("N18").Value = -83.661202185792

# ---------------------------------------------------------------------
# Row 19 — Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 35
$ws.Range("E19").Value = -25.714285714285
$ws.Range("F19").Value = 116
$ws.Range("G19").Value = 148
$ws.Range("H19").Value = -21.621621621621
$ws.Range("I19").Value = 1171
$ws.Range("J19").Value = 1398
$ws.Range("K19").Value = -16.237482117310
$ws.Range("L19").Value = -27.267080745341
$ws.Range("M19").Value = -24.108878807517
$ws.Range("N19").Value = -82.072872014696

# ---------------------------------------------------------------------
# Row 20 — G.L.A.
# ---------------------------------------------------------------------
Set-NumberCell "C20" 1 "J20"
Set-PlaceholderCell "D20" "0" "D23"
Set-PlaceholderCell "E20" "***.*" "N22"
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -71.428571428571
$ws.Range("I20").Value = 16
$ws.Range("K20").Value = -60
$ws.Range("L20").Value = -66.666666666666
$ws.Range("M20").Value = 6.666666666666
$ws.Range("N20").Value = -93.388429752066

# ---------------------------------------------------------------------
# Row 21 — TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 60
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = 15.384615384615
$ws.Range("F21").Value = 235
$ws.Range("G21").Value = 258
$ws.Range("H21").Value = -8.914728682170
$ws.Range("I21").Value = 2144
$ws.Range("J21").Value = 2376
$ws.Range("K21").Value = -9.764309764309
$ws.Range("L21").Value = -18.695487296169
$ws.Range("M21").Value = 5.304518664047
$ws.Range("N21").Value = -80.317635178555

# ---------------------------------------------------------------------
# Row 22 — Transit
# ---------------------------------------------------------------------
Set-NumberCell "C22" 2 "J22"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 134
$ws.Range("J22").Value = 119
$ws.Range("K22").Value = 12.605042016806
$ws.Range("L22").Value = -7.586206896551
$ws.Range("M22").Value = 41.052631578947

# ---------------------------------------------------------------------
# Row 24 — Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 75
$ws.Range("D24").Value = 76
$ws.Range("E24").Value = -1.315789473684
$ws.Range("F24").Value = 332
$ws.Range("G24").Value = 322
$ws.Range("H24").Value = 3.105590062111
$ws.Range("I24").Value = 2787
$ws.Range("J24").Value = 3002
$ws.Range("K24").Value = -7.161892071952
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -13.096351730589

# ---------------------------------------------------------------------
# Row 25 — Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 68
$ws.Range("D25").Value = 73
$ws.Range("E25").Value = -6.849315068493
$ws.Range("F25").Value = 285
$ws.Range("G25").Value = 279
$ws.Range("H25").Value = 2.150537634408
$ws.Range("I25").Value = 2389
$ws.Range("J25").Value = 2656
$ws.Range("K25").Value = -10.052710843373
$ws.Range("L25").Value = -5.907837731390

# ---------------------------------------------------------------------
# Row 26 — Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 28
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 95
$ws.Range("G26").Value = 87
$ws.Range("H26").Value = 9.195402298850
$ws.Range("I26").Value = 765
$ws.Range("J26").Value = 719
$ws.Range("K26").Value = 6.397774687065
$ws.Range("L26").Value = 2.960969044414
$ws.Range("M26").Value = 93.181818181818

# ---------------------------------------------------------------------
# Row 27 — UCR Rape*
# ---------------------------------------------------------------------
Set-NumberCell "C27" 2 "J27"
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 32
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 68.421052631578

# ---------------------------------------------------------------------
# Row 28 — Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 266.666666666667
$ws.Range("F28").Value = 63
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = 231.578947368421
$ws.Range("I28").Value = 277
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 73.125
$ws.Range("L28").Value = 77.564102564102

# ---------------------------------------------------------------------
# Rows 29-30 — no incidents this week: wipe the lone report-count figure
# back to the "0" placeholder.
# ---------------------------------------------------------------------
Set-PlaceholderCell "F29" "0" "G29"
Set-PlaceholderCell "F30" "0" "G30"

# ---------------------------------------------------------------------
# Row 31 — Hate Crimes (first figures reported in two years; fill in
# the previously-blank week-to-date / 28-day columns)
# ---------------------------------------------------------------------
Set-NumberCell "D31" 1 "I31"
Set-NumberCell "E31" -100 "K29"
Set-NumberCell "F31" 1 "I31"
Set-NumberCell "G31" 1 "I31"
Set-NumberCell "H31" 0 "K29"
$ws.Range("I31").Value = 13
$ws.Range("J31").Value = 19
$ws.Range("K31").Value = -31.578947368421
$ws.Range("L31").Value = 8.333333333333
